$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds 5 "test case" blocks, each describing a help-panel scenario
# for a different functionality. Blocks live at fixed row offsets
# (6-12, 15-21, 24-30, 33-39, 42-48); only the B/D cells on the two "Steps"
# rows (first two rows of the "Steps" table) carry functionality-specific
# text. Version 1.3 swaps the content of block 1 <-> block 3 (Avaliacoes <->
# Periodos Avaliativos) and block 2 <-> block 4 (Niveis das Competencias <->
# Competencias (portfolio)); block 5 (Perfis de Competencias) is unchanged.

$avaliacoes = @{
    B1 = "Lider de Pessoas acessa a funcionalidade de gestao de Avaliacoes a partir do menu inicial"
    D1 = "SYSTEM exibe a listagem das Avaliacoes cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
    B2 = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Avaliacoes"
}

$niveis = @{
    B1 = "Lider de Pessoas acessa a funcionalidade de gestao de Niveis das Competencias a partir do menu inicial"
    D1 = "SYSTEM exibe a listagem dos Niveis das Competencias cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
    B2 = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Niveis das Competencias"
}

$periodos = @{
    B1 = "Lider de Pessoas acessa a funcionalidade de gestao de Periodos Avaliativos a partir do menu inicial"
    D1 = "SYSTEM exibe a listagem dos Periodos Avaliativos cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
    B2 = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Periodos Avaliativos"
}

$competencias = @{
    B1 = "Lider de Pessoas acessa a funcionalidade de gestao de Competencias (portfolio) a partir do menu inicial"
    D1 = "SYSTEM exibe a listagem das Competencias (portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
    B2 = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Competencias (portfolio)"
}

function Set-Block($firstStepRow, $content) {
    $secondStepRow = $firstStepRow + 1
    $ws.Range("B$firstStepRow").Value = $content.B1
    $ws.Range("D$firstStepRow").Value = $content.D1
    $ws.Range("B$secondStepRow").Value = $content.B2
}

# Block 1 (rows 6-12): was Avaliacoes -> becomes Periodos Avaliativos
Set-Block 10 $periodos
# Block 2 (rows 15-21): was Niveis -> becomes Competencias (portfolio)
Set-Block 19 $competencias
# Block 3 (rows 24-30): was Periodos Avaliativos -> becomes Avaliacoes
Set-Block 28 $avaliacoes
# Block 4 (rows 33-39): was Competencias (portfolio) -> becomes Niveis das Competencias
Set-Block 37 $niveis
# Block 5 (rows 42-48, Perfis de Competencias) is left untouched.
